$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.625.86"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.801.75"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'227.64"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'32.85"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "'11.20"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "1.796.76"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "34.601.36"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").Value = "'68.95"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "0.0₃0806"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'247.70"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").Value = "'168.69"
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("D27").Value = "'16.61"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'4.09"
$ws.Range("E30").Value = "  +10.45%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0527"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "'3.83"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.24"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").Value = "1.434.44"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").Value = "'2.61"
$ws.Range("E36").Value = "  +7.99%  "
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "'1.07"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "'85.35"
$ws.Range("E40").Value = "  +6.10%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.944"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.41"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("D43").Value = "'2.77"
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").Value = "'13.92"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("E45").Value = "  +3.47%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'106.27"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.13%  "
